$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New block property rows to append (A..I), following the existing pattern.
$names = @("Grass1", "Grass2", "Grass3", "Grass4", "Grass5", "Crack1", "Crack2", "Crack3", "Crack4", "Crack5", "Treasure1 ", "Treasure2")

$row = 11
foreach ($name in $names) {
    $ws.Cells.Item($row, 2).Value = "string"
    $ws.Cells.Item($row, 2).NumberFormat = "@"

    $ws.Cells.Item($row, 3).Value = $false
    $ws.Cells.Item($row, 4).Value = $false
    $ws.Cells.Item($row, 5).Value = $false
    $ws.Cells.Item($row, 6).Value = $true

    $ws.Cells.Item($row, 7).Value = 0
    $ws.Cells.Item($row, 8).Value = 0

    $ws.Cells.Item($row, 9).Value = "Friend"
    $ws.Cells.Item($row, 9).NumberFormat = "@"

    $ws.Cells.Item($row, 1).Value = $name
    if ($row -le 20) {
        $ws.Cells.Item($row, 1).NumberFormat = "@"
    }

    $row = $row + 1
}

[void]$ws.Range("H24").Select()

Write-Host "Done"
